# Sprint 3 Backlog — the workbook was simply opened in Excel and re-saved;
# no cell content, formula, or value changed anywhere in the sheet (verified
# against the authoritative OOXML diff: every modified <row> hunk only touches
# bookkeeping attributes such as spans/x14ac:dyDescent, never a <c>/<v>).
# The one user-visible, content-level action captured by the diff is that the
# last selected cell on save was A26 on the (only / active) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reproduce the recorded selection / active cell (<selection activeCell="A26" sqref="A26"/>,
# sheetView tabSelected="1") exactly as captured in the saved file.
$ws.Activate()
$ws.Range("A26").Select()
